# Insert a new weekly record row before the existing row 154 (pushing the
# old rows 154-176 down to 155-177), then populate the new row 154 with the
# latest data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 154:176 down to 155:177 by inserting a new blank row at 154.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new observation.
$ws.Range("A154").Value = 4
$ws.Range("B154").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C154").Value = "Los Lagos"
$newDate = Get-Date -Year 2021 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("D154").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D154").Value = $newDate
$ws.Range("E154").Value = 10
$ws.Range("F154").Value = 100112037
$ws.Range("G154").Value = "Cebollín"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 180
$ws.Range("K154").Value = 5000
$ws.Range("L154").Value = 5500
$ws.Range("M154").Value = 5250
$ws.Range("N154").Value = "$/paquete 36 unidades"
$ws.Range("O154").Value = "Región Metropolitana"
$ws.Range("P154").Value = 146
$ws.Range("Q154").Value = 36
$ws.Range("R154").Value = "Hortaliza"
